$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DELIVER")

# ---------------------------------------------------------------------------
# The edit inserts a new "LAST SPRINTS" labelled section (mirroring the
# existing "DELIVERIES" section pattern) right before the existing
# IMAGE0 / description rows. Concretely (matching the structural diff):
#   - insert 2 new rows at row 43      (old 43 -> 45, old 44 -> 46, ...)
#   - insert 1 new row at row 47       (old 45 "DELIVERIES" header -> 48, ...)
# ---------------------------------------------------------------------------

$ws.Rows("43:44").Insert() | Out-Null
$ws.Rows("47:47").Insert() | Out-Null

# ---------------------------------------------------------------------------
# Row 46 currently still carries the old "thick bottom" row formatting that
# used to close off the IMAGE0/description block (it shifted down with the
# content). That visual divider now belongs to the new blank row 43 (before
# the new header) and the new blank row 47 (after the description, before
# the DELIVERIES header) instead.
# ---------------------------------------------------------------------------

# Copy the (now shifted) divider formatting onto the two new blank rows.
$ws.Range("A46").Copy() | Out-Null
$ws.Range("A43").PasteSpecial(-4122) | Out-Null
$ws.Range("A46").Copy() | Out-Null
$ws.Range("A47").PasteSpecial(-4122) | Out-Null

# Row 46 (the description row) no longer sits at the boundary, so drop the
# thick-bottom emphasis back to a plain row like row 45 above it.
$ws.Range("A45").Copy() | Out-Null
$ws.Range("A46").PasteSpecial(-4122) | Out-Null

$ws.Application.CutCopyMode = $false

# Restore the actual cell values after the format-only pastes.
$ws.Range("A43").Value = ""
$ws.Range("A47").Value = ""

# ---------------------------------------------------------------------------
# Build the new "LAST SPRINTS" section header row (row 44), copying the
# look of the existing "DELIVERIES" header (now at row 48) and swapping in
# the new shared string.
# ---------------------------------------------------------------------------

$ws.Range("A48:B48").Copy() | Out-Null
$ws.Range("A44:B44").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

$ws.Range("A44:B44").Merge() | Out-Null
$ws.Range("A44").Value = "LAST SPRINTS"

# ---------------------------------------------------------------------------
# The Table104 list object needs to track the data block that moved down by
# three rows (A54:F55 -> A57:F58).
# ---------------------------------------------------------------------------

$lo = $ws.ListObjects.Item("Table104")
$lo.Resize($ws.Range("A57:F58")) | Out-Null

$wb.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# The selection / scroll position shifted down along with the new content
# (the author had scrolled to and selected a cell a few rows further down).
# ---------------------------------------------------------------------------

$ws.Range("A46").Select() | Out-Null
